# Germany Regionalliga Bayern - update of league bases (28-06-2024 19:47)
# The edit re-orders the match-result rows for three groups of fixtures
# that share the same date (their ids/results got corrected), while
# leaving columns A (row index), C (Div) and D (Date) untouched.
#
# Row groups affected:
#   - rows 32 <-> 33            (simple swap)
#   - rows 42 -> 43 -> 44 -> 42  (3-way rotation)
#   - rows 233 <-> 235           (simple swap, row 234 untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that change: B, and E through AD (2, and 5..30)
$colIndexes = @(2) + @(5..30)

function Get-RowValues($row) {
    $vals = @{}
    foreach ($c in $colIndexes) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    foreach ($c in $colIndexes) {
        $ws.Cells.Item($row, $c).Value2 = $vals[$c]
    }
}

# --- Swap rows 32 and 33 ---
$row32 = Get-RowValues 32
$row33 = Get-RowValues 33
Set-RowValues 32 $row33
Set-RowValues 33 $row32

# --- Rotate rows 42 -> 43 -> 44 -> 42 ---
# new row42 = old row43 ; new row43 = old row44 ; new row44 = old row42
$row42 = Get-RowValues 42
$row43 = Get-RowValues 43
$row44 = Get-RowValues 44
Set-RowValues 42 $row43
Set-RowValues 43 $row44
Set-RowValues 44 $row42

# --- Swap rows 233 and 235 ---
$row233 = Get-RowValues 233
$row235 = Get-RowValues 235
Set-RowValues 233 $row235
Set-RowValues 235 $row233
